$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''304.77'
$ws.Range("E2").Value = '''3.40%'
$ws.Range("G2").Value = '''2'

$ws.Range("D3").Value = '''35.33'
$ws.Range("E3").Value = '''13.26%'
$ws.Range("G3").Value = '''2'

$ws.Range("D4").Value = '''5.182'
$ws.Range("E4").Value = '''5.03%'
$ws.Range("G4").Value = '''2'

$ws.Range("D5").Value = '''0.07807'
$ws.Range("E5").Value = '''5.35%'
$ws.Range("G5").Value = '''2'

$ws.Range("D6").Value = '''2.384'
$ws.Range("E6").Value = '''10.88%'
$ws.Range("G6").Value = '''2'

$ws.Range("D7").Value = '''8.052'
$ws.Range("E7").Value = '''4.11%'
$ws.Range("G7").Value = '''2'

$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = '''0.9365'
$ws.Range("E8").Value = '''2.41%'
$ws.Range("G8").Value = '''2'

$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = '''0.09811'
$ws.Range("E9").Value = '''12.23%'
$ws.Range("G9").Value = '''2'

$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '''0.1791'
$ws.Range("E10").Value = '''5.54%'
$ws.Range("G10").Value = '''2'

$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '''0.08573'
$ws.Range("E11").Value = '''2.67%'
$ws.Range("G11").Value = '''2'

$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '''0.03319'
$ws.Range("E12").Value = '''5.23%'
$ws.Range("G12").Value = '''2'

$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '''0.09913'
$ws.Range("E13").Value = '''-1.63%'
$ws.Range("G13").Value = '''2'

$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '''0.001502'
$ws.Range("E14").Value = '''-0.33%'
$ws.Range("G14").Value = '''2'

$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '''0.005846'
$ws.Range("E15").Value = '''-0.06%'
$ws.Range("G15").Value = '''2'

$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '''3.470'
$ws.Range("E16").Value = '''-0.83%'
$ws.Range("G16").Value = '''2'

$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = '''3.927'
$ws.Range("E17").Value = '''4.81%'
$ws.Range("G17").Value = '''2'

$ws.Range("D18").Value = '''2.137'
$ws.Range("E18").Value = '''2.93%'
$ws.Range("G18").Value = '''2'

$ws.Range("E19").Value = '''1.15%'
$ws.Range("G19").Value = '''2'

$ws.Range("E20").Value = '''2.89%'
$ws.Range("G20").Value = '''2'

$ws.Range("D21").Value = '''4.351'
$ws.Range("E21").Value = '''9.53%'
$ws.Range("G21").Value = '''2'

$ws.Range("D22").Value = '''0.2302'
$ws.Range("E22").Value = '''9.66%'
$ws.Range("G22").Value = '''2'

$ws.Range("D23").Value = '''0.04636'
$ws.Range("E23").Value = '''1.80%'
$ws.Range("G23").Value = '''2'

$ws.Range("D24").Value = '''0.001210'
$ws.Range("E24").Value = '''-0.11%'
$ws.Range("G24").Value = '''2'

$ws.Range("D25").Value = '''0.004393'
$ws.Range("E25").Value = '''-5.04%'
$ws.Range("G25").Value = '''2'

$ws.Range("D26").Value = '''0.0001300'
$ws.Range("E26").Value = '''-0.01%'
$ws.Range("G26").Value = '''2'

$ws.Range("D27").Value = '''0.0003394'
$ws.Range("E27").Value = '''-0.01%'
$ws.Range("G27").Value = '''2'

$ws.Range("G28").Value = '''2'

$ws.Range("G29").Value = '''2'

$ws.Range("G30").Value = '''2'

$ws.Range("G31").Value = '''2'

$ws.Range("G32").Value = '''2'

$ws.Range("G33").Value = '''2'

$ws.Range("G34").Value = '''2'

$ws.Range("G35").Value = '''2'

$ws.Range("G36").Value = '''2'

$ws.Range("G37").Value = '''2'

$ws.Range("G38").Value = '''2'

$ws.Range("D39").Value = '''0.01814'
$ws.Range("E39").Value = '''12.45%'
$ws.Range("G39").Value = '''2'

$ws.Range("D40").Value = '''0.04832'
$ws.Range("E40").Value = '''8.13%'
$ws.Range("G40").Value = '''2'

$ws.Range("D41").Value = '''0.007801'
$ws.Range("E41").Value = '''6.24%'
$ws.Range("G41").Value = '''2'

$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = '''0.1414'
$ws.Range("E42").Value = '''6.27%'
$ws.Range("G42").Value = '''2'

$ws.Range("B43").Value = 'Dexo'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QkL_pl546+dexo-dexo'
$ws.Range("D43").Value = '''0.008788'
$ws.Range("E43").Value = '''-1.97%'
$ws.Range("G43").Value = '''2'

$ws.Range("D44").Value = '''0.002078'
$ws.Range("E44").Value = '''6.64%'
$ws.Range("G44").Value = '''2'

$ws.Range("D45").Value = '''0.01006'
$ws.Range("E45").Value = '''6.66%'
$ws.Range("G45").Value = '''2'

$ws.Range("D46").Value = '''0.00006109'
$ws.Range("E46").Value = '''0.82%'
$ws.Range("G46").Value = '''2'

$ws.Range("E47").Value = '''-0.01%'
$ws.Range("G47").Value = '''2'

$ws.Range("D48").Value = '''2.779'
$ws.Range("E48").Value = '''19.55%'
$ws.Range("G48").Value = '''2'

$ws.Range("E49").Value = '''-31.01%'
$ws.Range("G49").Value = '''2'

$ws.Range("D50").Value = '''0.00002100'
$ws.Range("E50").Value = '''-0.01%'
$ws.Range("G50").Value = '''2'

$ws.Range("D51").Value = '''0.0002000'
$ws.Range("E51").Value = '''-0.01%'
$ws.Range("G51").Value = '''2'
